# End of Day 6/9
# Mark the unit tests for the exploration_check class (rows 15,16,18,19,20,21
# of the "Unittest necessary/complete?" column F) as Complete, switching their
# highlight from red ("Necessary") to green ("Complete") to match the rest of
# the sheet's convention, then move the view/selection toward the next class
# (mergeConstraints) to work on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$completeGreen = 5287936   # RGB(0,176,80) -> matches the existing "Complete" fill used elsewhere on the sheet

$rows = @(15, 16, 18, 19, 20, 21)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = "Complete"
    $cell.Interior.Color = $completeGreen
}

# Update the frozen-pane scroll position and active selection to reflect
# where work continues next (heading into the merge_constraints section).
$ws.Activate()
$ws.Range("A11").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A11").Select()
$ws.Range("D30").Select()
